$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New smoothened output values (column B), replacing the previous smoother results
$values = @(
    -9.5,
    -9.3333333333333304,
    -2.6666666666666599,
    0,
    26,
    28.3333333333333,
    65,
    82.6666666666666,
    124,
    115.333333333333,
    109.666666666666,
    136,
    162,
    188,
    207.666666666666,
    244.333333333333,
    287,
    303,
    356,
    413.33333333333297,
    477.33333333333297,
    505.33333333333297,
    527.33333333333303,
    577,
    635.66666666666595,
    683,
    716,
    751.33333333333303,
    806.66666666666595,
    883.66666666666595,
    971.33333333333303,
    1062.3333333333301,
    1125,
    1187.6666666666599,
    1231,
    1300,
    1355.3333333333301,
    1426.3333333333301,
    1503.6666666666599,
    1584.3333333333301,
    1668.6666666666599,
    1755.3333333333301,
    1856,
    1967,
    2030.3333333333301,
    2098.3333333333298,
    2196,
    2298,
    2416,
    2512.6666666666601,
    2618.3333333333298,
    2728.6666666666601,
    2831.6666666666601,
    2930.6666666666601,
    3010.3333333333298,
    3126,
    3242.3333333333298,
    3378,
    3485.3333333333298,
    3625.3333333333298,
    3752.3333333333298,
    3874.6666666666601,
    3992,
    4118.3333333333303,
    4250.6666666666597,
    4367.6666666666597,
    4505,
    4621.6666666666597,
    4753.3333333333303,
    4887.3333333333303,
    5025,
    5170.3333333333303,
    5318.6666666666597,
    5465.3333333333303,
    5614.6666666666597,
    5758.6666666666597,
    5934.6666666666597,
    6085,
    6252.6666666666597,
    6409.6666666666597,
    6598,
    6758.6666666666597,
    6924.3333333333303,
    7096,
    7241.6666666666597,
    7415.6666666666597,
    7566,
    7743.3333333333303,
    7895.3333333333303,
    8096.3333333333303,
    8300.3333333333303,
    8510.3333333333303,
    8696.3333333333303,
    8871.6666666666606,
    9042.3333333333303,
    9223.6666666666606,
    9429.3333333333303,
    9622,
    9820.3333333333303,
    9902.5
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $values[$i]
}
